$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Copy the style from E1 (bold header) to the new header cells
$ws.Range("A1").Copy() | Out-Null
$ws.Range("F1:H1").PasteSpecial(-4122) | Out-Null

# Boolean values for the Outliers_MAD columns (rows 2-12)
$values = @(
    @(0,0,0),
    @(0,0,0),
    @(1,1,1),
    @(0,0,1),
    @(0,0,0),
    @(0,0,0),
    @(1,0,1),
    @(0,0,1),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0)
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = [bool]$values[$i][0]
    $ws.Cells.Item($row, 7).Value = [bool]$values[$i][1]
    $ws.Cells.Item($row, 8).Value = [bool]$values[$i][2]
}
